# Insert a new price record as row 380 on the single data sheet.
# This pushes the existing rows 380-489 down to 381-490 (last row becomes 490),
# matching the diff: a new "Morrón rojo" record (fecha 45027) is inserted,
# and the sheet dimension grows from A1:R489 to A1:R490.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 380..489 down by one row.
$ws.Rows("380").Insert()

# Populate the newly inserted row 380 with the new record.
$ws.Cells.Item(380, 1).Value = 11
$ws.Cells.Item(380, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(380, 3).Value = "Bíobío"
$ws.Cells.Item(380, 4).Value = 45027
$ws.Cells.Item(380, 5).Value = 8
$ws.Cells.Item(380, 6).Value = 100112002
$ws.Cells.Item(380, 7).Value = "Pimiento"
$ws.Cells.Item(380, 8).Value = "Morrón rojo"
$ws.Cells.Item(380, 9).Value = "Primera"
$ws.Cells.Item(380, 10).Value = 220
$ws.Cells.Item(380, 11).Value = 9500
$ws.Cells.Item(380, 12).Value = 10000
$ws.Cells.Item(380, 13).Value = 9773
$ws.Cells.Item(380, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(380, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(380, 16).Value = 543
$ws.Cells.Item(380, 17).Value = 18
$ws.Cells.Item(380, 18).Value = "Hortaliza"
